# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets: row 2 (F2) 511 -> 513, row 3 (F3) 445 -> 447.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 513
    $ws.Range("F3").Value = 447
}
